# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the freshly scraped data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 171
    4  = 127
    5  = 1273
    6  = 17844
    7  = 344
    8  = 251
    10 = 6738
    12 = 153
    15 = 58
    16 = 6
    19 = 187
    23 = 29
    25 = 264
    26 = 967
    27 = 108
    28 = 5144
    29 = 531
    30 = 65
    31 = 11933
    32 = 1268
    34 = 197
    35 = 263
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 171
    4  = 127
    5  = 1273
    6  = 17844
    7  = 344
    8  = 251
    10 = 6738
    12 = 153
    15 = 58
    16 = 6
    19 = 187
    23 = 29
    25 = 264
    26 = 967
    27 = 108
    28 = 5144
    29 = 531
    32 = 65
    33 = 11933
    34 = 1268
    36 = 197
    37 = 263
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
